# "Add component error solved"
#
# Sheet1 had a stray, completely empty column O sitting between the
# Booking table (N) and the Attendance/PaySlip/Component/Plan/
# PlanComponent/Floor lookup tables that used to start at P. Remove the
# empty column so those tables shift left and start at O again (P->O,
# Q->P, R->Q, S->R, T->S, U->T), matching the rest of the sheet layout.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Columns("O:O").Delete()

# Tidy the widths of the now-shifted component/plan/floor columns so they
# fit their content again.
$ws.Range("R1:T19").EntireColumn.AutoFit()

# Restore the selection to where the author left off after the cleanup.
$ws.Range("S13").Select()
